$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.014.91"
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.591.74"
$ws.Range("E3").Value = "  +0.50%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.39"
$ws.Range("E5").Value = "  +0.49%  "

$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.480"
$ws.Range("E7").Value = "  +0.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.246"
$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.00"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0810"
$ws.Range("E11").Value = "  +2.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.815.99"
$ws.Range("E12").Value = "  +0.64%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.595.73"
$ws.Range("E13").Value = "  +0.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.00"
$ws.Range("E14").Value = "  -0.75%  "

$ws.Range("E15").Value = "  +1.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.988.87"
$ws.Range("E16").Value = "  +0.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.94"
$ws.Range("E17").Value = "  +1.96%  "

$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.36"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.24"
$ws.Range("E21").Value = "  +1.32%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.20"
$ws.Range("E22").Value = "  -1.51%  "

$ws.Range("E23").Value = "  +0.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.93"
$ws.Range("E24").Value = "  +14.14%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.68"
$ws.Range("E25").Value = "  +1.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("E27").Value = "  -7.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.16"
$ws.Range("E28").Value = "  +0.67%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.48"
$ws.Range("E29").Value = "  +0.53%  "

$ws.Range("E30").Value = "  +0.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0474"
$ws.Range("E31").Value = "  +0.79%  "

$ws.Range("E32").Value = "  -0.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.89"
$ws.Range("E33").Value = "  -4.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.47"
$ws.Range("E34").Value = "  -1.44%  "

$ws.Range("E35").Value = "  -0.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.126.36"
$ws.Range("E36").Value = "  +2.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0163"
$ws.Range("E37").Value = "  +8.29%  "

$ws.Range("E38").Value = "  +0.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.792"
$ws.Range("E39").Value = "  +2.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.32"
$ws.Range("E40").Value = "  -1.28%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.490"
$ws.Range("E41").Value = "  -2.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.778"
$ws.Range("E42").Value = "  -3.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.13"
$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.727.35"
$ws.Range("E44").Value = "  +0.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.17"
$ws.Range("E45").Value = "  -1.59%  "

$ws.Range("E46").Value = "  -0.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "53.51"
$ws.Range("E47").Value = "  +0.92%  "

$ws.Range("E48").Value = "  -1.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₇0969"
$ws.Range("E49").Value = "  -13.21%  "

$ws.Range("E50").Value = "  +0.08%  "

$ws.Range("E51").Value = "  +0.28%  "
